# Updates the cryptocurrency rows in the worksheet (prices, 1h volume %,
# and a couple of row reorderings) to match the refreshed data feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in columns B/C/D hold plain text that can look like numbers
# (e.g. "1.00", "26.78"); a leading apostrophe forces Excel to keep
# them as text instead of auto-converting to a numeric value.
$textPrefix = "'"

$ws.Range("D2").Value = $textPrefix + '29.531.60'
$ws.Range("E2").Value = '  +2.18%  '
$ws.Range("D3").Value = $textPrefix + '1.595.82'
$ws.Range("E3").Value = '  +1.03%  '
$ws.Range("E4").Value = '  +0.54%  '
$ws.Range("D5").Value = $textPrefix + '211.89'
$ws.Range("E5").Value = '  +0.10%  '
$ws.Range("E6").Value = '  -1.29%  '
$ws.Range("E7").Value = '  +0.58%  '
$ws.Range("D8").Value = $textPrefix + '26.78'
$ws.Range("E8").Value = '  +5.65%  '
$ws.Range("D9").Value = $textPrefix + '43.57'
$ws.Range("E9").Value = '  -2.86%  '
$ws.Range("D10").Value = $textPrefix + '0.250'
$ws.Range("E10").Value = '  +0.94%  '
$ws.Range("E11").Value = '  +0.91%  '
$ws.Range("E12").Value = '  +0.97%  '
$ws.Range("D13").Value = $textPrefix + '1.824.55'
$ws.Range("E13").Value = '  +1.12%  '
$ws.Range("D14").Value = $textPrefix + '1.595.20'
$ws.Range("E14").Value = '  +1.47%  '
$ws.Range("D15").Value = $textPrefix + '29.537.40'
$ws.Range("E15").Value = '  +2.25%  '
$ws.Range("E16").Value = '  +2.78%  '
$ws.Range("E17").Value = '  +1.02%  '
$ws.Range("D18").Value = $textPrefix + '63.73'
$ws.Range("E18").Value = '  +2.64%  '
$ws.Range("D19").Value = $textPrefix + '240.60'
$ws.Range("E19").Value = '  +3.74%  '
$ws.Range("D20").Value = $textPrefix + '7.54'
$ws.Range("E20").Value = '  +1.46%  '
$ws.Range("D21").Value = $textPrefix + '0.0₃0690'
$ws.Range("E21").Value = '  -0.29%  '
$ws.Range("E22").Value = '  +0.59%  '
$ws.Range("D23").Value = $textPrefix + '3.97'
$ws.Range("E23").Value = '  -0.20%  '
$ws.Range("D24").Value = $textPrefix + '9.20'
$ws.Range("E24").Value = '  +0.59%  '
$ws.Range("E25").Value = '  +0.27%  '
$ws.Range("D26").Value = $textPrefix + '154.89'
$ws.Range("E26").Value = '  +1.60%  '
$ws.Range("D27").Value = $textPrefix + '15.30'
$ws.Range("E27").Value = '  +2.29%  '
$ws.Range("E28").Value = '  +0.25%  '
$ws.Range("D29").Value = $textPrefix + '6.38'
$ws.Range("E29").Value = '  +0.84%  '
$ws.Range("E30").Value = '  +0.52%  '
$ws.Range("D31").Value = $textPrefix + '0.0476'
$ws.Range("E31").Value = '  +2.56%  '
$ws.Range("E32").Value = '  +0.12%  '
$ws.Range("D33").Value = $textPrefix + '3.21'
$ws.Range("E33").Value = '  +0.25%  '
# --- row 34/35: Maker and InternetComputer(DFINITY) swap places ---
$ws.Range("B34").Value = $textPrefix + 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = $textPrefix + 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").Value = $textPrefix + '3.13'
$ws.Range("E34").Value = '  +3.40%  '
$ws.Range("B35").Value = $textPrefix + 'Maker'
$ws.Range("C35").Value = $textPrefix + 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D35").Value = $textPrefix + '1.430.09'
$ws.Range("E35").Value = '  +0.70%  '
$ws.Range("E36").Value = '  +2.11%  '
$ws.Range("E37").Value = '  -2.07%  '
$ws.Range("D38").Value = $textPrefix + '2.82'
$ws.Range("E38").Value = '  +2.82%  '
$ws.Range("E39").Value = '  +0.56%  '
$ws.Range("E40").Value = '  +1.63%  '
$ws.Range("E41").Value = '  +2.82%  '
$ws.Range("D42").Value = $textPrefix + '1.95'
$ws.Range("E42").Value = '  +0.93%  '
$ws.Range("E43").Value = '  +6.58%  '
$ws.Range("D44").Value = $textPrefix + '53.40'
$ws.Range("E44").Value = '  +23.68%  '
# --- row 45/46: ARBITRUM and PaxDollar swap places ---
$ws.Range("B45").Value = $textPrefix + 'PaxDollar'
$ws.Range("C45").Value = $textPrefix + 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").Value = $textPrefix + '1.00'
$ws.Range("E45").Value = '  +0.42%  '
$ws.Range("B46").Value = $textPrefix + 'ARBITRUM'
$ws.Range("C46").Value = $textPrefix + 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D46").Value = $textPrefix + '0.797'
$ws.Range("E46").Value = '  +1.78%  '
$ws.Range("E47").Value = '  +16.62%  '
$ws.Range("D48").Value = $textPrefix + '65.39'
$ws.Range("E48").Value = '  +1.48%  '
$ws.Range("E49").Value = '  +0.18%  '
$ws.Range("D50").Value = $textPrefix + '1.737.13'
$ws.Range("E50").Value = '  +1.17%  '
$ws.Range("D51").Value = $textPrefix + '85.94'
$ws.Range("E51").Value = '  +0.78%  '
